$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.89 = 41642.53 pesos`n✅ 41642.53 pesos = 9.84 = 964.77 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 101.1
$ws2.Range("O10").Value = 4210.06
$ws2.Range("N12").Value = 4230
$ws2.Range("O12").Value = 98
